# Add a new "2020" data column (column R) to the right of the existing
# "2019" column (Q), mirroring the formatting of the corresponding Q cells
# and filling in the new year's values row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "R4"  = 2020
    "R5"  = 5
    "R6"  = 3.5
    "R7"  = 1.8
    "R8"  = 24.4
    "R9"  = 7.2
    "R10" = 2.9
    "R11" = 7.4
    "R12" = 4
    "R13" = 3.2
    "R14" = 3.5
}

foreach ($row in 4..14) {
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")

    # Clone the formatting (number format, borders, font, fill, alignment)
    # from the matching "2019" cell onto the new "2020" cell.
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null

    # Now write the new value for 2020.
    $dstCell.Value = $newValues["R$row"]
}

# Update the selection to the newly added column, matching what was active
# after the edit.
$ws.Range("R4:R14").Select() | Out-Null
